$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 9.208741113779867
$ws.Cells.Item(2, 3).Value = 4.539354369108066
$ws.Cells.Item(2, 5).Value = 19.54081999720889
$ws.Cells.Item(2, 6).Value = 41.44308217402294
$ws.Cells.Item(2, 7).Value = 34.80040057236256
$ws.Cells.Item(2, 8).Value = 15.82221624725277
$ws.Cells.Item(2, 9).Value = 22.84200665700756
$ws.Cells.Item(2, 10).Value = 8.450590942829677
$ws.Cells.Item(2, 11).Value = 9.477224571857272
$ws.Cells.Item(2, 13).Value = 17.9252899035782
$ws.Cells.Item(2, 14).Value = 19.63762078740946
$ws.Cells.Item(3, 2).Value = 8.960847332564697
$ws.Cells.Item(3, 3).Value = 4.358208924005166
$ws.Cells.Item(3, 5).Value = 19.45028667546163
$ws.Cells.Item(3, 6).Value = 41.36656660365244
$ws.Cells.Item(3, 7).Value = 34.84918848737769
$ws.Cells.Item(3, 8).Value = 15.8698950994418
$ws.Cells.Item(3, 9).Value = 22.92395237437858
$ws.Cells.Item(3, 10).Value = 8.467105432882777
$ws.Cells.Item(3, 11).Value = 9.321605766843465
$ws.Cells.Item(3, 13).Value = 17.82427926332689
$ws.Cells.Item(3, 14).Value = 19.70183387900351
$ws.Cells.Item(4, 2).Value = 8.80706270628637
$ws.Cells.Item(4, 3).Value = 4.244093993349057
$ws.Cells.Item(4, 5).Value = 19.3983861489304
$ws.Cells.Item(4, 6).Value = 41.32976643811801
$ws.Cells.Item(4, 7).Value = 34.88957252525545
$ws.Cells.Item(4, 8).Value = 15.90173952600342
$ws.Cells.Item(4, 9).Value = 22.97851285037408
$ws.Cells.Item(4, 10).Value = 8.477749708025204
$ws.Cells.Item(4, 11).Value = 9.226718091354323
$ws.Cells.Item(4, 13).Value = 17.76540038192556
$ws.Cells.Item(4, 14).Value = 19.74306799528815
$ws.Cells.Item(5, 2).Value = 8.744098856646197
$ws.Cells.Item(5, 3).Value = 4.196940808209475
$ws.Cells.Item(5, 5).Value = 19.37817907057669
$ws.Cells.Item(5, 6).Value = 41.31733979157597
$ws.Cells.Item(5, 7).Value = 34.90864241716294
$ws.Cells.Item(5, 8).Value = 15.91536202595342
$ws.Cells.Item(5, 9).Value = 23.00181279474091
$ws.Cells.Item(5, 10).Value = 8.48221453344873
$ws.Cells.Item(5, 11).Value = 9.18826827945715
$ws.Cells.Item(5, 13).Value = 17.74221619107715
$ws.Cells.Item(5, 14).Value = 19.76032699821053
$ws.Cells.Item(6, 2).Value = 8.733628974189322
$ws.Cells.Item(6, 3).Value = 4.189074174472257
$ws.Cells.Item(6, 5).Value = 19.37488107538586
$ws.Cells.Item(6, 6).Value = 41.31543179717714
$ws.Cells.Item(6, 7).Value = 34.91196648779669
$ws.Cells.Item(6, 8).Value = 15.91766301293878
$ws.Cells.Item(6, 9).Value = 23.00574608058658
$ws.Cells.Item(6, 10).Value = 8.48296360797086
$ws.Cells.Item(6, 11).Value = 9.181898357865498
$ws.Cells.Item(6, 13).Value = 17.73841590333355
$ws.Cells.Item(6, 14).Value = 19.76322041009797
$ws.Cells.Item(7, 2).Value = 8.806214612290225
$ws.Cells.Item(7, 3).Value = 4.243460595278853
$ws.Cells.Item(7, 5).Value = 19.39810979213609
$ws.Cells.Item(7, 6).Value = 41.32958843235501
$ws.Cells.Item(7, 7).Value = 34.88981914140005
$ws.Cells.Item(7, 8).Value = 15.90192063027465
$ws.Cells.Item(7, 9).Value = 22.97882276688693
$ws.Cells.Item(7, 10).Value = 8.477809406698906
$ws.Cells.Item(7, 11).Value = 9.226198594876386
$ws.Cells.Item(7, 13).Value = 17.76508440981039
$ws.Cells.Item(7, 14).Value = 19.74329890904669
$ws.Cells.Item(8, 2).Value = 9.123653654221677
$ws.Cells.Item(8, 3).Value = 4.477538756953811
$ws.Cells.Item(8, 5).Value = 19.50884854178868
$ws.Cells.Item(8, 6).Value = 41.41459313328139
$ws.Cells.Item(8, 7).Value = 34.81505319044051
$ws.Cells.Item(8, 8).Value = 15.83812235718973
$ws.Cells.Item(8, 9).Value = 22.86937952241565
$ws.Cells.Item(8, 10).Value = 8.456180700915521
$ws.Cells.Item(8, 11).Value = 9.423460145687889
$ws.Cells.Item(8, 13).Value = 17.88982124944749
$ws.Cells.Item(8, 14).Value = 19.65938724744436
$ws.Cells.Item(9, 2).Value = 9.729474677484449
$ws.Cells.Item(9, 3).Value = 4.910681920191025
$ws.Cells.Item(9, 5).Value = 19.7545233199762
$ws.Cells.Item(9, 6).Value = 41.66153524734903
$ws.Cells.Item(9, 7).Value = 34.75153007960296
$ws.Cells.Item(9, 8).Value = 15.73342193016839
$ws.Cells.Item(9, 9).Value = 22.68851234817755
$ws.Cells.Item(9, 10).Value = 8.41775067951971
$ws.Cells.Item(9, 11).Value = 9.813327908322904
$ws.Cells.Item(9, 13).Value = 18.15845914953627
$ws.Cells.Item(9, 14).Value = 19.50911059409066
$ws.Cells.Item(10, 2).Value = 10.15925493030568
$ws.Cells.Item(10, 3).Value = 5.209735735662123
$ws.Cells.Item(10, 5).Value = 19.95133801525892
$ws.Cells.Item(10, 6).Value = 41.89102188651928
$ws.Cells.Item(10, 7).Value = 34.75590692640617
$ws.Cells.Item(10, 8).Value = 15.66896765642511
$ws.Cells.Item(10, 9).Value = 22.57629297279307
$ws.Cells.Item(10, 10).Value = 8.391919781371877
$ws.Cells.Item(10, 11).Value = 10.09869747488639
$ws.Cells.Item(10, 13).Value = 18.36920010574296
$ws.Cells.Item(10, 14).Value = 19.4073127150367
$ws.Cells.Item(11, 2).Value = 10.35048620873135
$ws.Cells.Item(11, 3).Value = 5.341044334916726
$ws.Cells.Item(11, 5).Value = 20.0441683884303
$ws.Cells.Item(11, 6).Value = 42.00562316982126
$ws.Cells.Item(11, 7).Value = 34.76903516532674
$ws.Cells.Item(11, 8).Value = 15.64235919902097
$ws.Cells.Item(11, 9).Value = 22.52974742462678
$ws.Cells.Item(11, 10).Value = 8.380685231932953
$ws.Cells.Item(11, 11).Value = 10.22769453608081
$ws.Cells.Item(11, 13).Value = 18.46768501551772
$ws.Cells.Item(11, 14).Value = 19.36285235667756
$ws.Cells.Item(12, 2).Value = 10.42221414717156
$ws.Cells.Item(12, 3).Value = 5.390047520725121
$ws.Cells.Item(12, 5).Value = 20.07977272714772
$ws.Cells.Item(12, 6).Value = 42.05046390399065
$ws.Cells.Item(12, 7).Value = 34.77560963273713
$ws.Cells.Item(12, 8).Value = 15.63267377037052
$ws.Cells.Item(12, 9).Value = 22.5127710759647
$ws.Cells.Item(12, 10).Value = 8.376504809736126
$ws.Cells.Item(12, 11).Value = 10.27637740063656
$ws.Cells.Item(12, 13).Value = 18.50532954641135
$ws.Cells.Item(12, 14).Value = 19.34628067576991
$ws.Cells.Item(13, 2).Value = 10.40679789385366
$ws.Cells.Item(13, 3).Value = 5.379526431537364
$ws.Cells.Item(13, 5).Value = 20.07208498704891
$ws.Cells.Item(13, 6).Value = 42.04074284566121
$ws.Cells.Item(13, 7).Value = 34.77412239824996
$ws.Cells.Item(13, 8).Value = 15.63474231824687
$ws.Cells.Item(13, 9).Value = 22.51639832337653
$ws.Cells.Item(13, 10).Value = 8.377401859154796
$ws.Cells.Item(13, 11).Value = 10.26590077486559
$ws.Cells.Item(13, 13).Value = 18.49720695753038
$ws.Cells.Item(13, 14).Value = 19.34983794018886
$ws.Cells.Item(14, 2).Value = 10.35640147148338
$ws.Cells.Item(14, 3).Value = 5.345090506322786
$ws.Cells.Item(14, 5).Value = 20.04708866368506
$ws.Cells.Item(14, 6).Value = 42.00928345515705
$ws.Cells.Item(14, 7).Value = 34.76954391928965
$ws.Cells.Item(14, 8).Value = 15.64155454094685
$ws.Cells.Item(14, 9).Value = 22.52833774402733
$ws.Cells.Item(14, 10).Value = 8.380339827858428
$ws.Cells.Item(14, 11).Value = 10.2317032882118
$ws.Cells.Item(14, 13).Value = 18.470775181481
$ws.Cells.Item(14, 14).Value = 19.36148370082364
$ws.Cells.Item(15, 2).Value = 10.32544068682306
$ws.Cells.Item(15, 3).Value = 5.323902581784243
$ws.Cells.Item(15, 5).Value = 20.03183580092718
$ws.Cells.Item(15, 6).Value = 41.99020093425446
$ws.Cells.Item(15, 7).Value = 34.76694825671337
$ws.Cells.Item(15, 8).Value = 15.64577811197589
$ws.Cells.Item(15, 9).Value = 22.53573561426135
$ws.Cells.Item(15, 10).Value = 8.38214902664402
$ws.Cells.Item(15, 11).Value = 10.21073340183563
$ws.Cells.Item(15, 13).Value = 18.45462981070119
$ws.Cells.Item(15, 14).Value = 19.36865147322854
$ws.Cells.Item(16, 2).Value = 10.14666512779028
$ws.Cells.Item(16, 3).Value = 5.201055630769074
$ws.Cells.Item(16, 5).Value = 19.94533577212823
$ws.Cells.Item(16, 6).Value = 41.88373580795891
$ws.Cells.Item(16, 7).Value = 34.75527324595217
$ws.Cells.Item(16, 8).Value = 15.67076120950222
$ws.Cells.Item(16, 9).Value = 22.57942561432285
$ws.Cells.Item(16, 10).Value = 8.392664339948954
$ws.Cells.Item(16, 11).Value = 10.09024660531714
$ws.Cells.Item(16, 13).Value = 18.36281435264678
$ws.Cells.Item(16, 14).Value = 19.41025538361035
$ws.Cells.Item(17, 2).Value = 10.03584311613893
$ws.Cells.Item(17, 3).Value = 5.124450615942678
$ws.Cells.Item(17, 5).Value = 19.89309962775187
$ws.Cells.Item(17, 6).Value = 41.82102072213709
$ws.Cells.Item(17, 7).Value = 34.75096514771555
$ws.Cells.Item(17, 8).Value = 15.68678261218748
$ws.Cells.Item(17, 9).Value = 22.60738279914509
$ws.Cells.Item(17, 10).Value = 8.399247068177106
$ws.Cells.Item(17, 11).Value = 10.01608748668355
$ws.Cells.Item(17, 13).Value = 18.30714071350267
$ws.Cells.Item(17, 14).Value = 19.43625046183931
$ws.Cells.Item(18, 2).Value = 9.971702458254136
$ws.Cells.Item(18, 3).Value = 5.079946002826508
$ws.Cells.Item(18, 5).Value = 19.8633666577354
$ws.Cells.Item(18, 6).Value = 41.78591094441851
$ws.Cells.Item(18, 7).Value = 34.74953553788632
$ws.Cells.Item(18, 8).Value = 15.69625294015901
$ws.Cells.Item(18, 9).Value = 22.62388682455013
$ws.Cells.Item(18, 10).Value = 8.403081873712761
$ws.Cells.Item(18, 11).Value = 9.973358771705465
$ws.Cells.Item(18, 13).Value = 18.27536670526214
$ws.Cells.Item(18, 14).Value = 19.45137616045157
$ws.Cells.Item(19, 2).Value = 9.949919346292853
$ws.Cells.Item(19, 3).Value = 5.064802657728015
$ws.Cells.Item(19, 5).Value = 19.85335383230479
$ws.Cells.Item(19, 6).Value = 41.774189343628
$ws.Cells.Item(19, 7).Value = 34.7492314673108
$ws.Cells.Item(19, 8).Value = 15.69950324933687
$ws.Cells.Item(19, 9).Value = 22.62954753490284
$ws.Cells.Item(19, 10).Value = 8.404388629827892
$ws.Cells.Item(19, 11).Value = 9.958880280979153
$ws.Cells.Item(19, 13).Value = 18.26465194816777
$ws.Cells.Item(19, 14).Value = 19.45652738798357
$ws.Cells.Item(20, 2).Value = 10.0476821218413
$ws.Cells.Item(20, 3).Value = 5.132651552086577
$ws.Cells.Item(20, 5).Value = 19.89862813893939
$ws.Cells.Item(20, 6).Value = 41.82759742421479
$ws.Cells.Item(20, 7).Value = 34.75131524092099
$ws.Cells.Item(20, 8).Value = 15.68505068550998
$ws.Cells.Item(20, 9).Value = 22.60436283491029
$ws.Cells.Item(20, 10).Value = 8.398541298729254
$ws.Cells.Item(20, 11).Value = 10.02398989888472
$ws.Cells.Item(20, 13).Value = 18.31304178721414
$ws.Cells.Item(20, 14).Value = 19.4334652403278
$ws.Cells.Item(21, 2).Value = 10.37122331473088
$ws.Cells.Item(21, 3).Value = 5.355225018446426
$ws.Cells.Item(21, 5).Value = 20.05441862210056
$ws.Cells.Item(21, 6).Value = 42.01848486097617
$ws.Cells.Item(21, 7).Value = 34.77084521781973
$ws.Cells.Item(21, 8).Value = 15.63954301992706
$ws.Cells.Item(21, 9).Value = 22.52481320663256
$ws.Cells.Item(21, 10).Value = 8.379474873530118
$ws.Cells.Item(21, 11).Value = 10.24175278165115
$ws.Cells.Item(21, 13).Value = 18.47852953420132
$ws.Cells.Item(21, 14).Value = 19.35805589029716
$ws.Cells.Item(22, 2).Value = 10.57864178542149
$ws.Cells.Item(22, 3).Value = 5.496475102203637
$ws.Cells.Item(22, 5).Value = 20.15885723566514
$ws.Cells.Item(22, 6).Value = 42.15164323013647
$ws.Cells.Item(22, 7).Value = 34.79295277338971
$ws.Cells.Item(22, 8).Value = 15.61207821953821
$ws.Cells.Item(22, 9).Value = 22.47660920100924
$ws.Cells.Item(22, 10).Value = 8.367444201948935
$ws.Cells.Item(22, 11).Value = 10.38308795186679
$ws.Cells.Item(22, 13).Value = 18.58871550044427
$ws.Cells.Item(22, 14).Value = 19.31031260969631
$ws.Cells.Item(23, 2).Value = 10.4683299299777
$ws.Cells.Item(23, 3).Value = 5.421484705148042
$ws.Cells.Item(23, 5).Value = 20.10288433617784
$ws.Cells.Item(23, 6).Value = 42.079813780747
$ws.Cells.Item(23, 7).Value = 34.78029852329599
$ws.Cells.Item(23, 8).Value = 15.62652814600699
$ws.Cells.Item(23, 9).Value = 22.50198954069694
$ws.Cells.Item(23, 10).Value = 8.373825936737191
$ws.Cells.Item(23, 11).Value = 10.30776011997516
$ws.Cells.Item(23, 13).Value = 18.52973014761177
$ws.Cells.Item(23, 14).Value = 19.33565349740253
$ws.Cells.Item(24, 2).Value = 10.04233103371008
$ws.Cells.Item(24, 3).Value = 5.1289453490533
$ws.Cells.Item(24, 5).Value = 19.89612776816736
$ws.Cells.Item(24, 6).Value = 41.82462114875607
$ws.Cells.Item(24, 7).Value = 34.75115370198523
$ws.Cells.Item(24, 8).Value = 15.68583288153311
$ws.Cells.Item(24, 9).Value = 22.6057268180925
$ws.Cells.Item(24, 10).Value = 8.398860220362629
$ws.Cells.Item(24, 11).Value = 10.02041750635351
$ws.Cells.Item(24, 13).Value = 18.31037318244548
$ws.Cells.Item(24, 14).Value = 19.43472387572434
$ws.Cells.Item(25, 2).Value = 9.567925763851438
$ws.Cells.Item(25, 3).Value = 4.796655643322743
$ws.Cells.Item(25, 5).Value = 19.68511473337296
$ws.Cells.Item(25, 6).Value = 41.58621925815565
$ws.Cells.Item(25, 7).Value = 34.75977277130725
$ws.Cells.Item(25, 8).Value = 15.75955844963332
$ws.Cells.Item(25, 9).Value = 22.73381932896885
$ws.Cells.Item(25, 10).Value = 8.427723118730652
$ws.Cells.Item(25, 11).Value = 9.70783885008019
$ws.Cells.Item(25, 13).Value = 18.08334095570566
$ws.Cells.Item(25, 14).Value = 19.54824577722588
